$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

foreach ($r in 2..5) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value = 45175
    }
}
